# B6-PowerPoint.pptx edit replay
#
# The authoring commit:
#   1) Re-applied / changed the deck's theme so that the slide-master's
#      theme (ppt/theme/theme2.xml, the one actually driving the slides)
#      switched from the "Integral" / "Red Violet" palette to the plain
#      "Office Theme" palette (what used to live in the orphaned
#      ppt/theme/theme1.xml, only ever linked from the Notes Master).
#   2) As a direct, automatic consequence of the new theme being applied,
#      the three native tables in the deck (slides 14-16) - which were
#      using the *default* table style tied to the old theme - picked up
#      the new theme's default "No Style, Table Grid"-ish table style
#      ({6FF261D0-4F8D-49A8-8E19-9A69C6A0A660}) instead of the deck's
#      custom table style ({B7D3503B-44DE-4391-B2D5-4D81A5AE34E4}).
#
# Reproduce both effects through the exposed PowerPoint object model.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Swap the live theme's colour palette from "Red Violet"/Integral to
#    the plain Office palette. Table.Style can't be set directly (the
#    host requires ApplyStyle), and likewise the theme's colour scheme
#    is edited in place, one RGB swatch at a time, through
#    Slide.ThemeColorScheme - this is the same palette that backs every
#    slide (and hence every table / shape) in the deck.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$colors = $slide1.ThemeColorScheme

# Order of ThemeColorScheme items: dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink.
# RGB() packs as 0x00BBGGRR, i.e. R + G*256 + B*65536.
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72

# ---------------------------------------------------------------------
# 2) Re-point every table on the deck at the new theme's default table
#    style GUID.
# ---------------------------------------------------------------------
$newTableStyle = "{6FF261D0-4F8D-49A8-8E19-9A69C6A0A660}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
